# Add reportng / custom listener support: the AddCustomerTest sheet no longer
# carries an "alerttext" result column (the alert-text assertions moved to
# the listener) nor a per-row "runmode" column, so drop the stale data.

$wb = $excel.ActiveWorkbook

$wsAddCustomer = $wb.Worksheets.Item(1)   # AddCustomerTest
$wsOpenAccount = $wb.Worksheets.Item(2)   # OpenAccountTest
$wsTestSuite   = $wb.Worksheets.Item(3)   # test_suite

# Drop the runmode column (E) entirely from AddCustomerTest.
$null = $wsAddCustomer.Range("E1:E5").Delete()

# Keep the "alerttext" header (D1) but clear the now-unused per-row values.
$null = $wsAddCustomer.Range("D2:D5").ClearContents()

# Selection on AddCustomerTest moves off the old B5 cell.
$null = $wsAddCustomer.Range("D1").Select()

# test_suite becomes the active/selected sheet (was AddCustomerTest).
$null = $wsTestSuite.Activate()
$null = $wsTestSuite.Range("B14").Select()
